# Replace the "Fournisseur_XX" company names in column B (rows 2-11) with a
# list of last names ("noms de famille"), and apply bold formatting to the
# new values so they stand out, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastNames = @("Leroy", "Martin", "Mercier", "Moreau", "Muller", "Petit", "Robert", "Roger", "Rousseau", "Stéphane")

for ($i = 0; $i -lt $lastNames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $lastNames[$i]
}

$ws.Range("B2:B11").Font.Bold = $true

# Restore the view to show column A and leave a single-cell selection,
# as in the saved workbook.
$null = $ws.Range("A1").Select()
$null = $ws.Range("C23").Select()
